$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add B3 value 0 (new FormNumber entry)
$ws.Range("B3").Value = 0

# Clear the stray empty formatted cell R2
$ws.Range("R2").Clear()

# Update the active selection to A3
$ws.Range("A3").Select()
